$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows before the existing row 373, shifting the old
# rows 373-379 down to 375-381 (formatting/styles carry down with them).
$ws.Rows("373:374").Insert()

# New row 373: Albahaca, Primera, Región Metropolitana, $/docena de matas
$ws.Cells.Item(373, 1).Value = 6
$ws.Cells.Item(373, 2).Value = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Cells.Item(373, 3).Value = "Metropolitana"
$ws.Cells.Item(373, 4).Value = 44595
$ws.Cells.Item(373, 5).Value = 13
$ws.Cells.Item(373, 6).Value = 100112052
$ws.Cells.Item(373, 7).Value = "Albahaca"
$ws.Cells.Item(373, 8).Value = "Sin especificar"
$ws.Cells.Item(373, 9).Value = "Primera"
$ws.Cells.Item(373, 10).Value = 1380
$ws.Cells.Item(373, 11).Value = 2500
$ws.Cells.Item(373, 12).Value = 3000
$ws.Cells.Item(373, 13).Value = 2703
$ws.Cells.Item(373, 14).Value = "`$/docena de matas"
$ws.Cells.Item(373, 15).Value = "Región Metropolitana"
$ws.Cells.Item(373, 16).Value = 450
$ws.Cells.Item(373, 17).Value = 6
$ws.Cells.Item(373, 18).Value = "Hortaliza"

# New row 374: Albahaca, Segunda, Región Metropolitana, $/docena de matas
$ws.Cells.Item(374, 1).Value = 6
$ws.Cells.Item(374, 2).Value = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Cells.Item(374, 3).Value = "Metropolitana"
$ws.Cells.Item(374, 4).Value = 44595
$ws.Cells.Item(374, 5).Value = 13
$ws.Cells.Item(374, 6).Value = 100112052
$ws.Cells.Item(374, 7).Value = "Albahaca"
$ws.Cells.Item(374, 8).Value = "Sin especificar"
$ws.Cells.Item(374, 9).Value = "Segunda"
$ws.Cells.Item(374, 10).Value = 490
$ws.Cells.Item(374, 11).Value = 2000
$ws.Cells.Item(374, 12).Value = 2000
$ws.Cells.Item(374, 13).Value = 2000
$ws.Cells.Item(374, 14).Value = "`$/docena de matas"
$ws.Cells.Item(374, 15).Value = "Región Metropolitana"
$ws.Cells.Item(374, 16).Value = 333
$ws.Cells.Item(374, 17).Value = 6
$ws.Cells.Item(374, 18).Value = "Hortaliza"
